# Applies weekly update of Alcachofa prices at Terminal Hortofrutícola Agro Chillán.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 updates ---
$ws.Range("D3").Value = 44426
$ws.Range("O3").Value = "Región del Maule"

# --- Row 4 updates ---
$ws.Range("D4").Value = 44420
$ws.Range("O4").Value = "Provincia del Elquí"

# --- Row 5 updates ---
$ws.Range("D5").Value = 44427
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 13500
$ws.Range("P5").Value = 338

# --- Row 6 updates ---
$ws.Range("D6").Value = 44417
$ws.Range("H6").Value = "Madrigal"
$ws.Range("I6").Value = "Primera"
$ws.Range("K6").Value = 15000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 388

# --- Row 7 updates ---
$ws.Range("D7").Value = 44432
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("P7").Value = 362

# --- Row 8 updates ---
$ws.Range("D8").Value = 44446
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 12500
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12750
$ws.Range("P8").Value = 319

# --- Row 9 updates ---
$ws.Range("D9").Value = 44435

# --- New row 10 (appended); match formatting/style of row 9 (D column uses the date style) ---
$ws.Range("D10").NumberFormat = $ws.Range("D9").NumberFormat

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44399
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Española"
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 15500
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15750
$ws.Range("N10").Value = "$/caja 40 unidades"
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 394
$ws.Range("Q10").Value = 40
$ws.Range("R10").Value = "Hortaliza"
